$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append: dates 2021-04-27 .. 2021-05-02
# (serials 44313..44318), with "nuovi pos.", "somma mobile 7gg.",
# "somma mobile 7gg. per 100mila abitanti"
$newRows = @(
    @{ Row = 239; A = 44313; B = 0; C = 1; D = 53.73455131649651 },
    @{ Row = 240; A = 44314; B = 0; C = 1; D = 53.73455131649651 },
    @{ Row = 241; A = 44315; B = 0; C = 0; D = 0 },
    @{ Row = 242; A = 44316; B = 0; C = 0; D = 0 },
    @{ Row = 243; A = 44317; B = 0; C = 0; D = 0 },
    @{ Row = 244; A = 44318; B = 0; C = 0; D = 0 }
)

# Template row used to copy the date cell formatting (column A) onto the
# newly appended rows, matching the style already used throughout the sheet.
$templateDateCell = $ws.Cells.Item(238, 1)
$templateDateCell.Copy()

foreach ($r in $newRows) {
    $rowIndex = $r.Row

    $ws.Cells.Item($rowIndex, 1).Value = $r.A
    $ws.Cells.Item($rowIndex, 1).PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($rowIndex, 2).Value = $r.B
    $ws.Cells.Item($rowIndex, 3).Value = $r.C
    $ws.Cells.Item($rowIndex, 4).Value = $r.D
}

$excel.CutCopyMode = $false
